# Added static wait in test case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChangePassword")

# Swap the old/new password test values on row 2
$ws.Range("A2").Value = "test1234"
$ws.Range("B2").Value = "test123"
$ws.Range("C2").Value = "test123"

# Move the active selection from D2 to D3
$ws.Activate()
$ws.Range("D3").Select()
